$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet plumbing: insert "2022-Q1" between "2021-Q4" and "总计".
#
#    Re-create "总计" so that the brand-new "2022-Q1" sheet can claim the
#    lowest free sheetId (2), leaving "总计" with sheetId 3 - matching how
#    the workbook looks after a human inserted a sheet in the middle in
#    real Excel.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$null = $zj.Delete()

$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

$zj2 = $wb.Worksheets.Add()
$zj2.Name = "总计"

# fix tab order back to 2021-Q4, 2022-Q1, 总计
$zjFix = $wb.Worksheets.Item("总计")
$q1Fix = $wb.Worksheets.Item("2022-Q1")
$zjFix.Move($null, $q1Fix)

$q4Fix = $wb.Worksheets.Item("2021-Q4")
$q1Fix2 = $wb.Worksheets.Item("2022-Q1")
$q4Fix.Move($q1Fix2)

# ---------------------------------------------------------------------------
# 2) Populate the new "2022-Q1" sheet (same layout as "2021-Q4").
# ---------------------------------------------------------------------------
$ref = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Item("2022-Q1")

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Cells.Item(2, 1).Value = 0
$ws.Range("B2").Value = "'011056"
$ws.Range("C2").Value = "博时汇兴回报一年持有期灵活配置混合"
$ws.Range("D2").Value = "'107.57"
$ws.Range("E2").Value = "'67.02"
$ws.Range("F2").Value = "'2.33"
$ws.Range("G2").Value = "'2.5064"
$ws.Cells.Item(2, 8).Value = 9

$ws.Cells.Item(3, 1).Value = 1
$ws.Range("B3").Value = "'011927"
$ws.Range("C3").Value = "博时汇誉回报灵活配置混合型证券投资基金A"
$ws.Range("D3").Value = "'1.30"
$ws.Range("E3").Value = "'68.12"
$ws.Range("F3").Value = "'5.71"
$ws.Range("G3").Value = "'0.0742"
$ws.Cells.Item(3, 8).Value = 2

$ws.Cells.Item(4, 1).Value = 2
$ws.Range("B4").Value = "'011928"
$ws.Range("C4").Value = "博时汇誉回报灵活配置混合型证券投资基金C"
$ws.Range("D4").Value = "'0.13"
$ws.Range("E4").Value = "'68.12"
$ws.Range("F4").Value = "'5.71"
$ws.Range("G4").Value = "'0.0074"
$ws.Cells.Item(4, 8).Value = 2

# replicate formatting from the reference sheet (header row + index column)
# so style indices line up (s="2" on headers/index, default style elsewhere)
# and any quote-prefix formatting picked up above is cleared again.
$ref.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ref.Range("A2:A4").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

$ref.Range("B2:H4").Copy()
$ws.Range("B2:H4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Populate the "总计" sheet: insert the 2022-Q1 row above 2021-Q4.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$zj.Range("B1").Value = "日期"
$zj.Range("C1").Value = "持有数量(只)"
$zj.Range("D1").Value = "持有市值(亿元)"

$zj.Cells.Item(2, 1).Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Cells.Item(2, 3).Value = 3
$zj.Cells.Item(2, 4).Value = 2.59

$zj.Cells.Item(3, 1).Value = 1
$zj.Range("B3").Value = "2021-Q4"
$zj.Cells.Item(3, 3).Value = 4
$zj.Cells.Item(3, 4).Value = 3.61

$ref.Range("B1:D1").Copy()
$zj.Range("B1:D1").PasteSpecial(-4122)

$ref.Range("A2:A3").Copy()
$zj.Range("A2:A3").PasteSpecial(-4122)
